$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J18").Value = 0.03
$ws.Range("K18").Value = 0.03
$ws.Range("L18").Value = 0.03
$ws.Range("M18").Value = 0.03
$ws.Range("Q18").Value = 0.035
$ws.Range("R18").Value = 0.17

$ws.Range("J19").Select()
